# Swap the "Top Sender(s)/Max Incoming Interactions" (B:C) pair with the
# "Top Receiver(s)/Max Outgoing Interactions" (D:E) pair for every data row
# (rows 2-16), leaving the header row (row 1) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2

    $ws.Cells.Item($r, 2).Value = $dVal
    $ws.Cells.Item($r, 3).Value = $eVal
    $ws.Cells.Item($r, 4).Value = $bVal
    $ws.Cells.Item($r, 5).Value = $cVal
}
